# "Generate Report for Handoff"
#
# This CI job regenerates the localization-status report. For the
# d3aed0d9-7321-49bf-b45a-1381e51fbd82 file, fresh handoff/handback
# timestamps were produced, which updates:
#   - Overview!D7            (Latest Handoff Date for that file)
#   - zh-cn!E7                (Latest Handoff Datetime for that file)
#   - de-de!E7                (Latest Handoff Datetime for that file)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-03-24 12:52:53"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-24 12:52:48"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-24 12:52:53"
